$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 45142, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '2a amarillo', 340, 12000, 13000, 12500, '$/caja 20 kilos', 'Región de Coquimbo', 625, 20),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 45142, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sutil De Gase', 'Primera', 300, 29000, 30000, 29500, '$/caja 24 kilos', 'Perú', 1229, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 45142, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 300, 34000, 35000, 34500, '$/caja 24 kilos', 'Perú', 1438, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44469, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '2a amarillo', 250, 10000, 11000, 10500, '$/caja 20 kilos', 'Región Metropolitana', 525, 20),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44910, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '2a amarillo', 300, 21000, 22000, 21500, '$/caja 20 kilos', 'Región Metropolitana', 1075, 20),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44910, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '3a amarillo', 300, 17000, 18000, 17500, '$/caja 20 kilos', 'Región Metropolitana', 875, 20),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44802, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sutil De Gase', 'Primera', 200, 46000, 47000, 46500, '$/caja 24 kilos', 'Perú', 1938, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44753, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 250, 17000, 18000, 17500, '$/caja 18 kilos', 'Perú', 972, 18),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44260, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sutil De Gase', 'Primera', 200, 27000, 28000, 27500, '$/caja 24 kilos', 'Perú', 1146, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44260, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 300, 22000, 23000, 22500, '$/caja 24 kilos', 'Perú', 938, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44258, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '2a amarillo', 300, 27000, 28000, 27500, '$/caja 20 kilos', 'Región de Coquimbo', 1375, 20),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44692, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '2a amarillo', 250, 19000, 20000, 19500, '$/caja 20 kilos', 'Región de Coquimbo', 975, 20),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44491, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sutil De Gase', 'Primera', 200, 55000, 56000, 55500, '$/caja 24 kilos', 'Perú', 2312, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44491, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 200, 45000, 46000, 45500, '$/caja 24 kilos', 'Perú', 1896, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44407, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sutil De Gase', 'Primera', 250, 31000, 32000, 31500, '$/caja 24 kilos', 'Perú', 1312, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44860, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '1a amarillo', 500, 14000, 15000, 14400, '$/caja 20 kilos', 'Región de Arica y Parinacota', 720, 20),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44860, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 750, 36000, 37000, 36533, '$/caja 24 kilos', 'Perú', 1522, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44447, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '2a amarillo', 270, 10000, 11000, 10500, '$/caja 20 kilos', 'Región de Coquimbo', 525, 20),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44431, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sutil De Gase', 'Primera', 200, 30000, 32000, 31000, '$/caja 24 kilos', 'Perú', 1292, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44431, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 300, 31000, 32000, 31500, '$/caja 24 kilos', 'Perú', 1312, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44382, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sutil De Gase', 'Primera', 180, 32000, 33000, 32556, '$/caja 24 kilos', 'Perú', 1356, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44420, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '3a amarillo', 250, 10000, 11000, 10500, '$/caja 20 kilos', 'Región Metropolitana', 525, 20),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44270, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 250, 29000, 30000, 29500, '$/caja 24 kilos', 'Perú', 1229, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 45138, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sutil De Gase', 'Primera', 450, 31000, 32000, 31444, '$/caja 24 kilos', 'Perú', 1310, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 45138, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 200, 26000, 27000, 26750, '$/caja 24 kilos', 'Colombia', 1115, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 45138, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 300, 28000, 29000, 28500, '$/caja 24 kilos', 'Perú', 1188, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44952, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '2a amarillo', 500, 25000, 26000, 25600, '$/caja 20 kilos', 'Región de O''Higgins', 1280, 20),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44715, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 300, 32000, 33000, 32500, '$/caja 24 kilos', 'Perú', 1354, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44603, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 350, 35000, 36000, 35500, '$/caja 24 kilos', 'Perú', 1479, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 45140, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '1a amarillo', 300, 16000, 17000, 16500, '$/caja 20 kilos', 'Región de Valparaíso', 825, 20),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 45140, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 250, 28000, 29000, 28600, '$/caja 24 kilos', 'Perú', 1192, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44610, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 300, 35000, 36000, 35500, '$/caja 24 kilos', 'Perú', 1479, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 45119, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 260, 29000, 30000, 29481, '$/caja 24 kilos', 'Perú', 1228, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44624, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sutil De Gase', 'Primera', 200, 46000, 47000, 46500, '$/caja 24 kilos', 'Perú', 1938, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44624, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 300, 45000, 46000, 45500, '$/caja 24 kilos', 'Perú', 1896, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 45043, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '2a plateado', 250, 2700, 28000, 12820, '$/caja 20 kilos', 'Región de O''Higgins', 641, 20),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44526, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sutil De Gase', 'Primera', 200, 21000, 22000, 21500, '$/caja 24 kilos', 'Perú', 896, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44526, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 300, 26000, 27000, 26500, '$/caja 24 kilos', 'Perú', 1104, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44341, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '3a plateado', 250, 10000, 11000, 10500, '$/caja 20 kilos', 'Región Metropolitana', 525, 20),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44524, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '2a amarillo', 250, 14000, 15000, 14500, '$/caja 20 kilos', 'Región de Coquimbo', 725, 20),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 45068, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sutil De Gase', 'Segunda', 150, 26000, 27000, 26333, '$/caja 24 kilos', 'Perú', 1097, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44631, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sutil De Gase', 'Primera', 250, 39000, 40000, 39500, '$/caja 24 kilos', 'Perú', 1646, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44631, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 300, 44000, 45000, 44500, '$/caja 24 kilos', 'Perú', 1854, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44645, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 300, 36000, 37000, 36500, '$/caja 24 kilos', 'Perú', 1521, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44846, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '2a amarillo', 300, 12000, 13000, 12500, '$/caja 20 kilos', 'Región de Coquimbo', 625, 20),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44846, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sutil De Gase', 'Primera', 250, 36000, 37000, 36500, '$/caja 24 kilos', 'Perú', 1521, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44846, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 250, 30000, 31000, 30500, '$/caja 24 kilos', 'Perú', 1271, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44237, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '2a amarillo', 250, 26000, 27000, 26500, '$/caja 20 kilos', 'Región de Coquimbo', 1325, 20),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44237, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sutil De Gase', 'Primera', 200, 21000, 22000, 21500, '$/caja 24 kilos', 'Perú', 896, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44237, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 200, 22000, 23000, 22500, '$/caja 24 kilos', 'Perú', 938, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44953, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 550, 29000, 31000, 29909, '$/caja 24 kilos', 'Perú', 1246, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44417, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sutil De Gase', 'Primera', 300, 32000, 33000, 32500, '$/caja 24 kilos', 'Perú', 1354, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44417, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 300, 30000, 31000, 30500, '$/caja 24 kilos', 'Perú', 1271, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44979, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '2a plateado', 300, 30000, 32000, 31000, '$/caja 20 kilos', 'Región de Coquimbo', 1550, 20),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 44300, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sin especificar', '2a amarillo', 270, 20000, 21000, 20500, '$/caja 20 kilos', 'Región de Coquimbo', 1025, 20),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 45133, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Sutil De Gase', 'Primera', 160, 31000, 32000, 31500, '$/caja 24 kilos', 'Perú', 1312, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 45133, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 200, 26000, 27000, 26500, '$/caja 24 kilos', 'Colombia', 1104, 24),
  @(1, 'Agrícola del Norte S.A. de Arica', 'Arica y Parinacota', 45133, 15, 'Fruta', 100102, 'Cítricos', 100102003, 'Limón', 'Tahití', 'Primera', 200, 28000, 29000, 28500, '$/caja 24 kilos', 'Perú', 1188, 24)
)

$nrows = $rows.Count
$ncols = $rows[0].Count
$data = New-Object 'object[,]' $nrows,$ncols
for ($i = 0; $i -lt $nrows; $i++) {
  for ($j = 0; $j -lt $ncols; $j++) {
    $data[$i,$j] = $rows[$i][$j]
  }
}

$startRow = 419
$endRow = $startRow + $nrows - 1
$ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 20)).Value = $data

# Ensure the date column (D) on the brand-new rows carries the same date
# number format as the rest of column D (existing rows keep their style).
$ws.Range("D474:D476").NumberFormat = "YYYY-MM-DD HH:MM:SS"
